# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (A) used four emoji as status glyphs:
#   📕 -> -3      (rouge / red)
#   📘 -> ⚠️      (bleu / blue)
#   📙 -> +3      (orange)
#   📗 -> ✅      (vert / green)
#
# Replace every occurrence in column A (rows 2..last) with the new plain
# text / emoji markers. "-3" and "+3" look numeric to Excel, so force the
# cell to Text format before writing them or Excel would silently store
# -3 as a number (and turn "+3" into the number 3, dropping the sign).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cur = $cell.Value2

    if ($cur -eq "📕") {
        $cell.NumberFormat = "@"
        $cell.Value2 = "-3"
    }
    elseif ($cur -eq "📘") {
        $cell.Value2 = "⚠️"
    }
    elseif ($cur -eq "📙") {
        $cell.NumberFormat = "@"
        $cell.Value2 = "+3"
    }
    elseif ($cur -eq "📗") {
        $cell.Value2 = "✅"
    }
}
